# Fix excel import ncc
# - Rename header D1 "Số lượng quy đổi" -> "Tỷ lệ quy đổi"
# - Rename header F1 "Đơn vị tính" -> "Đơn vị mua"
# - Set explicit column widths for A, D, E, F
# - Move active selection to H8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Tỷ lệ quy đổi"
$ws.Range("F1").Value = "Đơn vị mua"

# ColumnWidth is a character-width property that Excel's COM layer snaps to
# a 1/6-character (~pixel) grid before storing the OOXML `width` attribute
# (stored = (Round(ColumnWidth*6)+5)/6). The values below are chosen so the
# stored width lands on (or as close as possible to) the target widths
# 15 / 13.28515625 / 10.7109375 / 12 for columns A / D / E / F.
$ws.Columns.Item(1).ColumnWidth = 14.166666666666666
$ws.Columns.Item(4).ColumnWidth = 12.5
$ws.Columns.Item(5).ColumnWidth = 9.833333333333334
$ws.Columns.Item(6).ColumnWidth = 11.166666666666666

$ws.Range("H8").Select() | Out-Null
